$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly update swaps the data of row 2 with row 4, and row 3 with row 5
# (dates, variety, quality, volume, prices and unit of sale).

# --- Row 2 <-> Row 4 ---
$ws.Range("D2").Value = 44505
$ws.Range("K2").Value = "Californiana(o)"
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 100
$ws.Range("Q2").Value = "$/bandeja 10 kilos"

$ws.Range("D4").Value = 44902
$ws.Range("K4").Value = "Golden Nugget"
$ws.Range("L4").Value = "Especial"
$ws.Range("M4").Value = 60
$ws.Range("Q4").Value = "$/caja 10 kilos"

# --- Row 3 <-> Row 5 ---
$ws.Range("D3").Value = 44505
$ws.Range("M3").Value = 50
$ws.Range("N3").Value = 15000
$ws.Range("O3").Value = 15000
$ws.Range("P3").Value = 15000
$ws.Range("Q3").Value = "$/bandeja 10 kilos"
$ws.Range("S3").Value = 1500

$ws.Range("D5").Value = 44902
$ws.Range("M5").Value = 70
$ws.Range("N5").Value = 13000
$ws.Range("O5").Value = 13000
$ws.Range("P5").Value = 13000
$ws.Range("Q5").Value = "$/caja 10 kilos"
$ws.Range("S5").Value = 1300
